$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Sam"
$ws.Cells.Item(5, 3).Value = 3423

# Row 6
$ws.Cells.Item(6, 2).Value = "Sarah"

# Row 7
$ws.Cells.Item(7, 1).Value = "i"
$ws.Cells.Item(7, 2).Value = "Debby"
$ws.Cells.Item(7, 3).Value = "Row"

# Row 8
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 3).Value = "Blackmore"

# Row 9
$ws.Cells.Item(9, 1).Value = 9

# Column widths (bestFit recompute for the new, longer data)
$ws.Columns.Item(2).ColumnWidth = 4.9
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666

# Update selection to match diff (A1:C9 selection with active cell C9)
$ws.Range("A1:C9").Select() | Out-Null
$ws.Range("C9").Activate() | Out-Null
